# Applies the scheduled-runner profit recalculation update to the Leve
# profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Only the derived
# price/profit columns (H:N) on specific rows are refreshed; columns A:G
# (name/item/level/exp/gil/amount/item id) are left untouched.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value2 = 10152.321
$ws.Range("I62").Value2 = 15177.8
$ws.Range("J62").Value2 = 4353.6924
$ws.Range("K62").Value2 = 15177.8
$ws.Range("L62").Value2 = 4353.6924
$ws.Range("M62").Value2 = -14553.8
$ws.Range("N62").Value2 = -5601.6924
$ws.Range("H65").Value2 = 10152.321
$ws.Range("I65").Value2 = 15177.8
$ws.Range("J65").Value2 = 4353.6924
$ws.Range("K65").Value2 = 75889
$ws.Range("L65").Value2 = 21768.462
$ws.Range("M65").Value2 = -72769
$ws.Range("N65").Value2 = -28008.462
$ws.Range("H116").Value2 = 3521.5715
$ws.Range("I116").Value2 = 3562.875
$ws.Range("J116").Value2 = 3466.5
$ws.Range("K116").Value2 = 3562.875
$ws.Range("L116").Value2 = 3466.5
$ws.Range("M116").Value2 = -120.875
$ws.Range("N116").Value2 = -10350.5
$ws.Range("H132").Value2 = 2071.527
$ws.Range("I132").Value2 = 1156.6119
$ws.Range("K132").Value2 = 3469.835700000001
$ws.Range("M132").Value2 = -939.8357000000005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 2128.2144
$ws.Range("I2").Value2 = 2568.5
$ws.Range("J2").Value2 = 1798
$ws.Range("K2").Value2 = 2568.5
$ws.Range("L2").Value2 = 1798
$ws.Range("M2").Value2 = -2455.5
$ws.Range("N2").Value2 = -2024
$ws.Range("H45").Value2 = 1110.8
$ws.Range("I45").Value2 = 711.1667
$ws.Range("J45").Value2 = 1710.25
$ws.Range("K45").Value2 = 711.1667
$ws.Range("L45").Value2 = 1710.25
$ws.Range("M45").Value2 = -334.1667
$ws.Range("N45").Value2 = -2464.25
$ws.Range("H61").Value2 = 364223.2
$ws.Range("I61").Value2 = 269760.03
$ws.Range("K61").Value2 = 269760.03
$ws.Range("M61").Value2 = -269548.03
$ws.Range("H63").Value2 = 2966.6667
$ws.Range("H66").Value2 = 2966.6667
$ws.Range("H110").Value2 = 585.1667
$ws.Range("I110").Value2 = 562.2
$ws.Range("J110").Value2 = 700
$ws.Range("K110").Value2 = 562.2
$ws.Range("L110").Value2 = 700
$ws.Range("M110").Value2 = 1482.8
$ws.Range("N110").Value2 = -4790
$ws.Range("H116").Value2 = 2128.2144
$ws.Range("I116").Value2 = 2568.5
$ws.Range("J116").Value2 = 1798
$ws.Range("K116").Value2 = 2568.5
$ws.Range("L116").Value2 = 1798
$ws.Range("M116").Value2 = -274.5
$ws.Range("N116").Value2 = -6386
$ws.Range("H122").Value2 = 2933.0527
$ws.Range("I122").Value2 = 3044.1904
$ws.Range("J122").Value2 = 2795.7646
$ws.Range("K122").Value2 = 9132.5712
$ws.Range("L122").Value2 = 8387.293799999999
$ws.Range("M122").Value2 = -6682.5712
$ws.Range("N122").Value2 = -13287.2938
$ws.Range("H136").Value2 = 364223.2
$ws.Range("I136").Value2 = 269760.03
$ws.Range("K136").Value2 = 809280.0900000001
$ws.Range("M136").Value2 = -806730.0900000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 2128.2144
$ws.Range("I3").Value2 = 2568.5
$ws.Range("J3").Value2 = 1798
$ws.Range("K3").Value2 = 2568.5
$ws.Range("L3").Value2 = 1798
$ws.Range("M3").Value2 = -2454.5
$ws.Range("N3").Value2 = -2026

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 818.5
$ws.Range("I16").Value2 = 835.2143
$ws.Range("J16").Value2 = 760
$ws.Range("K16").Value2 = 835.2143
$ws.Range("L16").Value2 = 760
$ws.Range("M16").Value2 = -548.2143
$ws.Range("N16").Value2 = -1334
$ws.Range("H62").Value2 = 1427561.4
$ws.Range("I62").Value2 = 2417970.2
$ws.Range("J62").Value2 = 3848.5625
$ws.Range("K62").Value2 = 2417970.2
$ws.Range("L62").Value2 = 3848.5625
$ws.Range("M62").Value2 = -2417346.2
$ws.Range("N62").Value2 = -5096.5625
$ws.Range("H63").Value2 = 29000
$ws.Range("J63").Value2 = 29000
$ws.Range("L63").Value2 = 29000
$ws.Range("N63").Value2 = -30372
$ws.Range("H65").Value2 = 1427561.4
$ws.Range("I65").Value2 = 2417970.2
$ws.Range("J65").Value2 = 3848.5625
$ws.Range("K65").Value2 = 12089851
$ws.Range("L65").Value2 = 19242.8125
$ws.Range("M65").Value2 = -12086731
$ws.Range("N65").Value2 = -25482.8125
$ws.Range("H66").Value2 = 29000
$ws.Range("J66").Value2 = 29000
$ws.Range("L66").Value2 = 87000
$ws.Range("N66").Value2 = -93864
$ws.Range("H113").Value2 = 818.5
$ws.Range("I113").Value2 = 835.2143
$ws.Range("J113").Value2 = 760
$ws.Range("K113").Value2 = 835.2143
$ws.Range("L113").Value2 = 760
$ws.Range("M113").Value2 = 1334.7857
$ws.Range("N113").Value2 = -5100
$ws.Range("H134").Value2 = 1511.5454
$ws.Range("I134").Value2 = 864.7742
$ws.Range("K134").Value2 = 2594.3226
$ws.Range("M134").Value2 = -59.32259999999997

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value2 = 724
$ws.Range("I5").Value2 = 439.63635
$ws.Range("J5").Value2 = 1766.6666
$ws.Range("K5").Value2 = 1318.90905
$ws.Range("L5").Value2 = 5299.9998
$ws.Range("M5").Value2 = -1206.90905
$ws.Range("N5").Value2 = -5523.9998
$ws.Range("H34").Value2 = 593.05457
$ws.Range("J34").Value2 = 603.0741
$ws.Range("L34").Value2 = 1809.2223
$ws.Range("N34").Value2 = -1977.2223
$ws.Range("H39").Value2 = 3270.3
$ws.Range("J39").Value2 = 3900
$ws.Range("L39").Value2 = 11700
$ws.Range("N39").Value2 = -12288
$ws.Range("H113").Value2 = 12195707
$ws.Range("I113").Value2 = 16129612
$ws.Range("J113").Value2 = 600.2
$ws.Range("K113").Value2 = 48388836
$ws.Range("L113").Value2 = 1800.6
$ws.Range("M113").Value2 = -48386666
$ws.Range("N113").Value2 = -6140.6
$ws.Range("H132").Value2 = 4367.2593
$ws.Range("I132").Value2 = 2401.1428
$ws.Range("J132").Value2 = 6484.615
$ws.Range("K132").Value2 = 21610.2852
$ws.Range("L132").Value2 = 58361.535
$ws.Range("M132").Value2 = -19080.2852
$ws.Range("N132").Value2 = -63421.535
$ws.Range("H135").Value2 = 724
$ws.Range("I135").Value2 = 439.63635
$ws.Range("J135").Value2 = 1766.6666
$ws.Range("K135").Value2 = 3956.72715
$ws.Range("L135").Value2 = 15899.9994
$ws.Range("M135").Value2 = -1421.72715
$ws.Range("N135").Value2 = -20969.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value2 = 4308.8887
$ws.Range("I80").Value2 = 4761.7646
$ws.Range("J80").Value2 = 2909.0908
$ws.Range("K80").Value2 = 4761.7646
$ws.Range("L80").Value2 = 2909.0908
$ws.Range("M80").Value2 = -3763.7646
$ws.Range("N80").Value2 = -4905.0908
$ws.Range("H83").Value2 = 4308.8887
$ws.Range("I83").Value2 = 4761.7646
$ws.Range("J83").Value2 = 2909.0908
$ws.Range("K83").Value2 = 23808.823
$ws.Range("L83").Value2 = 14545.454
$ws.Range("M83").Value2 = -18816.823
$ws.Range("N83").Value2 = -24529.454
$ws.Range("H102").Value2 = 3734.7307
$ws.Range("I102").Value2 = 2102.1875
$ws.Range("K102").Value2 = 2102.1875
$ws.Range("M102").Value2 = -480.1875
$ws.Range("H126").Value2 = 3372.7778
$ws.Range("I126").Value2 = 3284.4443
$ws.Range("J126").Value2 = 3416.9443
$ws.Range("K126").Value2 = 9853.332900000001
$ws.Range("L126").Value2 = 10250.8329
$ws.Range("M126").Value2 = -7383.332900000001
$ws.Range("N126").Value2 = -15190.8329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value2 = 3739.8667
$ws.Range("I122").Value2 = 3233
$ws.Range("J122").Value2 = 4077.7778
$ws.Range("K122").Value2 = 9699
$ws.Range("L122").Value2 = 12233.3334
$ws.Range("M122").Value2 = -7249
$ws.Range("N122").Value2 = -17133.3334
$ws.Range("H132").Value2 = 6970.1904
$ws.Range("I132").Value2 = 2046.95
$ws.Range("K132").Value2 = 6140.85
$ws.Range("M132").Value2 = -3610.85

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value2 = 2106.2666
$ws.Range("I81").Value2 = 1133.3334
$ws.Range("J81").Value2 = 3565.6667
$ws.Range("K81").Value2 = 2266.6668
$ws.Range("L81").Value2 = 7131.3334
$ws.Range("M81").Value2 = -1205.6668
$ws.Range("N81").Value2 = -9253.3334
$ws.Range("H84").Value2 = 2106.2666
$ws.Range("I84").Value2 = 1133.3334
$ws.Range("J84").Value2 = 3565.6667
$ws.Range("K84").Value2 = 11333.334
$ws.Range("L84").Value2 = 35656.667
$ws.Range("M84").Value2 = -6029.333999999999
$ws.Range("N84").Value2 = -46264.667
$ws.Range("H107").Value2 = 1123.8334
$ws.Range("I107").Value2 = 988.6
$ws.Range("J107").Value2 = 1800
$ws.Range("K107").Value2 = 2965.8
$ws.Range("L107").Value2 = 5400
$ws.Range("M107").Value2 = -1045.8
$ws.Range("N107").Value2 = -9240
$ws.Range("H136").Value2 = 17433208
$ws.Range("I136").Value2 = 25026664
$ws.Range("K136").Value2 = 75079992
$ws.Range("M136").Value2 = -75077442
